$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.234.09'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').Value = '1.858.44'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('D4').Value = '''0.9997'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '''0.7114'
$ws.Range('E5').Value = '  +0.99%  '
$ws.Range('D6').Value = '''237.91'
$ws.Range('E6').Value = '  -0.43%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '''0.08157'
$ws.Range('E8').Value = '  +9.55%  '
$ws.Range('D9').Value = '''0.3042'
$ws.Range('E9').Value = '  -0.43%  '
$ws.Range('D10').Value = '''23.22'
$ws.Range('E10').Value = '  -1.02%  '
$ws.Range('D11').Value = '''0.08195'
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('D12').Value = '1.866.34'
$ws.Range('E12').Value = '  -1.46%  '
$ws.Range('D13').Value = '''5.173'
$ws.Range('E13').Value = '  -0.91%  '
$ws.Range('D14').Value = '''0.7082'
$ws.Range('E14').Value = '  -2.87%  '
$ws.Range('D15').Value = '''89.52'
$ws.Range('E15').Value = '  +0.32%  '
$ws.Range('D16').Value = '29.231.72'
$ws.Range('E16').Value = '  -0.15%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '''0.000007914'
$ws.Range('E17').Value = '  +3.32%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').Value = '''5.794'
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('D19').Value = '''13.36'
$ws.Range('E19').Value = '  +1.79%  '
$ws.Range('D20').Value = '''237.90'
$ws.Range('E20').Value = '  -0.47%  '
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').Value = '2.100.37'
$ws.Range('E22').Value = '  -2.98%  '
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').Value = '''7.414'
$ws.Range('E24').Value = '  -2.48%  '
$ws.Range('D25').Value = '''162.52'
$ws.Range('E25').Value = '  +0.82%  '
$ws.Range('D26').Value = '''0.1464'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').Value = '''8.964'
$ws.Range('E27').Value = '  -0.68%  '
$ws.Range('D28').Value = '''18.10'
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('D29').Value = '''1.959'
$ws.Range('E29').Value = '  -0.84%  '
$ws.Range('D30').Value = '''1.428'
$ws.Range('E30').Value = '  +1.22%  '
$ws.Range('D31').Value = '''4.407'
$ws.Range('E31').Value = '  -2.33%  '
$ws.Range('E32').Value = '  -1.02%  '
$ws.Range('D33').Value = '''4.024'
$ws.Range('E33').Value = '  +0.26%  '
$ws.Range('D34').Value = '''0.05219'
$ws.Range('E34').Value = '  +0.24%  '
$ws.Range('D35').Value = '''1.168'
$ws.Range('E35').Value = '  -1.86%  '
$ws.Range('D36').Value = '''0.7082'
$ws.Range('E36').Value = '  +0.14%  '
$ws.Range('D37').Value = '''0.9986'
$ws.Range('E37').Value = '  -3.96%  '
$ws.Range('D38').Value = '''2.672'
$ws.Range('E38').Value = '  +0.56%  '
$ws.Range('D39').Value = '''0.01861'
$ws.Range('E39').Value = '  -0.52%  '
$ws.Range('D40').Value = '''2.729'
$ws.Range('E40').Value = '  +1.82%  '
$ws.Range('D41').Value = '1.142.18'
$ws.Range('E41').Value = '  +7.00%  '
$ws.Range('D42').Value = '''0.9229'
$ws.Range('E42').Value = '  -2.10%  '
$ws.Range('D43').Value = '''0.4286'
$ws.Range('E43').Value = '  -0.56%  '
$ws.Range('D44').Value = '''5.876'
$ws.Range('E44').Value = '  -2.62%  '
$ws.Range('D45').Value = '''70.13'
$ws.Range('E45').Value = '  -0.73%  '
$ws.Range('D46').Value = '''0.9999'
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').Value = '''102.26'
$ws.Range('E47').Value = '  -1.31%  '
$ws.Range('D48').Value = '''1.776'
$ws.Range('E48').Value = '  +1.52%  '
$ws.Range('D49').Value = '2.002.63'
$ws.Range('E49').Value = '  -1.41%  '
$ws.Range('D50').Value = '''9.194'
$ws.Range('E50').Value = '  +0.97%  '
$ws.Range('D51').Value = '''6.971'
$ws.Range('E51').Value = '  -1.23%  '
